$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Add files via upload" re-save of the 109th Congress, session 2
# legislative-activity resume sheet. Corrects labeling/data errors left
# over from the prior upload:
#
# Relabel column A. The original sheet mislabeled the "Measures passed"
# block as bare "Senate bills" / "House bills" / ... rows, and used ad-hoc
# string codes ("S11857", "H9325", "E2257") in place of real numeric
# figures for "Pages of proceedings" / "Extension of remarks". This pass
# fixes the labels and converts those three cells to plain numbers.
$ws.Range("A2").Value = "Congress"
$ws.Range("A3").Value = "Session"
$ws.Range("A4").Value = "Start Date"
$ws.Range("A5").Value = "End Date"
$ws.Range("A6").Value = "Days in session"
$ws.Range("A7").Value = "Time in session"
$ws.Range("A8").Value = "Pages of proceedings"
$ws.Range("A9").Value = "Extensions of remarks"
$ws.Range("A10").Value = "Public bills enacted into law"
$ws.Range("A11").Value = "Private bills enacted into law"
$ws.Range("A12").Value = "Bills in conference"
$ws.Range("A13").Value = "Measures passed, total"
$ws.Range("A14").Value = "     Measures passed, Senate bills"
$ws.Range("A15").Value = "     Measures passed, House bills"
$ws.Range("A16").Value = "     Measures passed, Senate joint resolutions"
$ws.Range("A17").Value = "     Measures passed, House joint resolutions"
$ws.Range("A18").Value = "     Measures passed, Senate concurrent resolutions"
$ws.Range("A19").Value = "     Measures passed, House concurrent resolutions"
$ws.Range("A20").Value = "     Measures passed, Simple resolutions"
$ws.Range("A21").Value = "Measures reported, total"
$ws.Range("A22").Value = "     Measures reported, Senate bills"
$ws.Range("A23").Value = "     Measures reported, House bills"
$ws.Range("A24").Value = "     Measures reported, Senate joint resolutions"
$ws.Range("A25").Value = "     Measures reported, House joint resolutions"
$ws.Range("A26").Value = "     Measures reported, Senate concurrent resolutions"
$ws.Range("A27").Value = "     Measures reported, House concurrent resolutions"
$ws.Range("A28").Value = "     Measures reported, Simple resolutions"
$ws.Range("A29").Value = "Special reports"
$ws.Range("A30").Value = "Conference reports"
$ws.Range("A31").Value = "Measures pending on calendar"
$ws.Range("A32").Value = "Measures introduced, total"
$ws.Range("A33").Value = "     Measures introduced, Bills"
$ws.Range("A34").Value = "     Measures introduced, Joint resolutions"
$ws.Range("A35").Value = "     Measures introduced, Concurrent resolutions"
$ws.Range("A36").Value = "     Measures introduced, Simple resolutions"
$ws.Range("A37").Value = "Quorum calls"
$ws.Range("A38").Value = "Yea-and-nay votes"
$ws.Range("A39").Value = "Recorded votes"
$ws.Range("A40").Value = "Bills vetoed"
$ws.Range("A41").Value = "Vetoes overridden"

# "Time in session": normalize the missing space after the comma in the
# Senate figure ("1027 hrs,48'" -> "1027 hrs, 48'"); House figure unchanged.
$ws.Range("B7").Value = "1027 hrs, 48'"
$ws.Range("C7").Value = "850 hrs, 19'"

# "Pages of proceedings" / "Extension of remarks" used placeholder text
# codes instead of the real page-count numbers.
$ws.Range("B8").Value = 11857
$ws.Range("C8").Value = 9325
$ws.Range("C9").Value = 2257

# Leave the selection where the author left it when they saved.
[void]$ws.Range("B8").Select()
